$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.341936826705933
$ws.Range("B1").Value = 4.428929805755615
$ws.Range("C1").Value = 3.302167177200317
$ws.Range("D1").Value = 0.8947980403900146
$ws.Range("E1").Value = 0.4714532494544983
